{"js": "// Replace the old placeholder-style <id> values with the final downloaded\n// tc/tcn ids:\n//   <id>p092v_a1</id> -> <id>p092v_1</id>\n//   <id>p092v_a2</id> -> <id>p092v_2</id>\n//   <id>p092v_a3</id> -> <id>p092v_3</id>\nconst replacements = [\n  [\"<id>p092v_a1</id>\", \"<id>p092v_1</id>\"],\n  [\"<id>p092v_a2</id>\", \"<id>p092v_2</id>\"],\n  [\"<id>p092v_a3</id>\", \"<id>p092v_3</id>\"]\n];\n\nfor (const [oldText, newText] of replacements) {\n  const results = context.document.body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (const range of results.items) {\n    range.insertText(newText, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "# Replace the old placeholder-style <id> values with the final downloaded\n# tc/tcn ids:\n#   <id>p092v_a1</id> -> <id>p092v_1</id>\n#   <id>p092v_a2</id> -> <id>p092v_2</id>\n#   <id>p092v_a3</id> -> <id>p092v_3</id>\n$word = New-Object -ComObject Word.Application\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @{ Old = \"<id>p092v_a1</id>\"; New = \"<id>p092v_1</id>\" },\n    @{ Old = \"<id>p092v_a2</id>\"; New = \"<id>p092v_2</id>\" },\n    @{ Old = \"<id>p092v_a3</id>\"; New = \"<id>p092v_3</id>\" }\n)\n\nforeach ($pair in $replacements) {\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $pair.Old\n    $find.Replacement.Text = $pair.New\n    $find.Forward = $true\n    $find.Wrap = 1\n    $find.MatchCase = $true\n    $find.MatchWholeWord = $false\n    $find.MatchWildcards = $false\n    $find.Execute($find.Text, $find.MatchCase, $false, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2)\n}\n"}
